$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 32, shifting rows 32:128 down to 33:129
# (this reproduces the "new weekly record inserted, rest pushed down" edit).
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new record's data.
$ws.Range("A32").Value = 8
$ws.Range("B32").Value = "Terminal La Palmera de La Serena"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = (Get-Date -Year 2021 -Month 8 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 100112012
$ws.Range("G32").Value = "Espinaca"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 3400
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 500
$ws.Range("M32").Value = 450
$ws.Range("N32").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O32").Value = "Provincia del Elquí"
$ws.Range("P32").Value = 900
$ws.Range("Q32").Value = 0.5
$ws.Range("R32").Value = "Hortaliza"
